$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.495.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.06%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.507.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.72%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.22%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.21%  '

# Row 7
$ws.Range('E7').Value = '  +0.15%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.515'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.33%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.506.95'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.74%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.60%  '

# Row 11
$ws.Range('E11').Value = '  -0.51%  '

# Row 12
$ws.Range('E12').Value = '  -5.01%  '

# Row 13
$ws.Range('E13').Value = '  -2.70%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.971.31'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.53%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.233.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.32%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.35%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.61%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.515.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.41%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.84%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.35%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.41%  '

# Row 22
$ws.Range('E22').Value = '  -5.85%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.52%  '

# Row 24
$ws.Range('E24').Value = '  +0.07%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.63%  '

# Row 26
$ws.Range('E26').Value = '  -5.95%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.37%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.640.56'
$ws.Range('D28').Style = 'Normal'

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.994'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.55%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0902'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.55%  '

# Row 31
$ws.Range('E31').Value = '  -2.75%  '

# Row 32
$ws.Range('E32').Value = '  -2.85%  '

# Row 33
$ws.Range('E33').Value = '  +0.85%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.97%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.08%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.115'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.43%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '151.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.06%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.11%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.60'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.22%  '

# Row 40
$ws.Range('E40').Value = '  -0.02%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.46%  '

# Row 42
$ws.Range('E42').Value = '  -2.93%  '

# Row 43
$ws.Range('E43').Value = '  -5.84%  '

# Row 44
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -14.09%  '

# Row 45
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.08%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.45%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '143.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.59%  '

# Row 48
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.531'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.78%  '

# Row 49
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.20%  '

# Row 50
$ws.Range('E50').Value = '  -5.27%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.587'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.59%  '
